$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "'05/50/5000"
$ws.Cells.Item(2,2).Value = "'5000.00"
$ws.Cells.Item(2,3).Value = "'5000.00"
$ws.Cells.Item(2,4).Value = "'5000.00"
$ws.Cells.Item(2,5).Value = "'5000.00"
$ws.Cells.Item(2,6).Value = "'0.00"
$ws.Cells.Item(2,7).Value = "'100.00"

$ws.Range("A2:G2").Style = "Normal"
